$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrección columna posición de las clasificaciones
# The position values in column A (rows 2-7) were off by one;
# decrement each by 1 so the ranking starts at 0 instead of 1.
for ($r = 2; $r -le 7; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    $cell.Value = $current - 1
}
